$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.791376948356628
$ws.Range("B1").Value = 3.882137060165405
$ws.Range("C1").Value = 1.357315540313721
$ws.Range("D1").Value = 0.8500217199325562
$ws.Range("E1").Value = 0.4603198766708374
